# Update countries & provincias Spain
# - Chequia's row gets refreshed data; Australia keeps its old data but the
#   two countries swap rows 36/37 (Chequia now sorts above Australia).
# - A handful of other countries (Estados Unidos, Espana, Alemania) get
#   updated running totals.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Estados Unidos (row 4)
$ws.Range("B4").Value = 707718
$ws.Range("C4").Value = 30148
$ws.Range("E4").Value = 611631
$ws.Range("G4").Value = 2312
$ws.Range("H4").Value = 36929

# Espana (row 5)
$ws.Range("B5").Value = 190839
$ws.Range("C5").Value = 5891
$ws.Range("E5").Value = 96040
$ws.Range("G5").Value = 687
$ws.Range("H5").Value = 20002

# Alemania (row 8)
$ws.Range("B8").Value = 141397
$ws.Range("C8").Value = 3699
$ws.Range("E8").Value = 53931
$ws.Range("G8").Value = 300
$ws.Range("H8").Value = 4352

# Row 36 becomes Chequia with refreshed figures
$ws.Range("A36").Value = "Chequia"
$ws.Range("B36").Value = 6549
$ws.Range("C36").Value = 116
$ws.Range("D36").Value = 1174
$ws.Range("E36").Value = 5202
$ws.Range("F36").Value = 82
$ws.Range("G36").Value = 4
$ws.Range("H36").Value = 173

# Row 37 becomes Australia, carrying the figures Chequia vacated from row 36
$ws.Range("A37").Value = "Australia"
$ws.Range("B37").Value = 6526
$ws.Range("C37").Value = 58
$ws.Range("D37").Value = 3821
$ws.Range("E37").Value = 2640
$ws.Range("F37").Value = 60
$ws.Range("G37").Value = 2
$ws.Range("H37").Value = 65
